# Scheduled data refresh: update Leve profit calculations across all profession sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 458.7857
$ws.Range("J2").Value = 605.75
$ws.Range("L2").Value = 605.75
$ws.Range("N2").Value = -831.75
$ws.Range("H6").Value = 209.7
$ws.Range("I6").Value = 209.7
$ws.Range("K6").Value = 629.0999999999999
$ws.Range("M6").Value = -517.0999999999999
$ws.Range("H31").Value = 88.5
$ws.Range("I31").Value = 88.5
$ws.Range("K31").Value = 265.5
$ws.Range("M31").Value = -35.5
$ws.Range("H39").Value = 100.07692
$ws.Range("I39").Value = 104.25
$ws.Range("J39").Value = 50
$ws.Range("K39").Value = 312.75
$ws.Range("L39").Value = 150
$ws.Range("M39").Value = -16.75
$ws.Range("N39").Value = -742
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("H86").Value = 4882.75
$ws.Range("I86").Value = 2665.6667
$ws.Range("J86").Value = 6213
$ws.Range("K86").Value = 2665.6667
$ws.Range("L86").Value = 6213
$ws.Range("M86").Value = -1542.6667
$ws.Range("N86").Value = -8459
$ws.Range("H89").Value = 4882.75
$ws.Range("I89").Value = 2665.6667
$ws.Range("J89").Value = 6213
$ws.Range("K89").Value = 13328.3335
$ws.Range("L89").Value = 31065
$ws.Range("M89").Value = -7712.333500000001
$ws.Range("N89").Value = -42297
$ws.Range("H106").Value = 0
$ws.Range("I106").Value = 0
$ws.Range("K106").Value = 0
$ws.Range("M106").ClearContents()
$ws.Range("H116").Value = 0
$ws.Range("I116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("M116").ClearContents()
$ws.Range("H132").Value = 7600.6
$ws.Range("I132").Value = 4003
$ws.Range("K132").Value = 12009
$ws.Range("M132").Value = -9479

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H23").Value = 10000
$ws.Range("J23").Value = 10000
$ws.Range("L23").Value = 10000
$ws.Range("N23").Value = -10518
$ws.Range("H32").Value = 3806.7827
$ws.Range("I32").Value = 3328.05
$ws.Range("K32").Value = 3328.05
$ws.Range("M32").Value = -3041.05
$ws.Range("H37").Value = 7500
$ws.Range("H45").Value = 1852.3334
$ws.Range("I45").Value = 1065
$ws.Range("J45").Value = 2246
$ws.Range("K45").Value = 1065
$ws.Range("L45").Value = 2246
$ws.Range("M45").Value = -688
$ws.Range("N45").Value = -3000
$ws.Range("H74").Value = 1113.4615
$ws.Range("I74").Value = 862.8
$ws.Range("K74").Value = 862.8
$ws.Range("M74").Value = 11.20000000000005
$ws.Range("H77").Value = 1113.4615
$ws.Range("I77").Value = 862.8
$ws.Range("K77").Value = 4314
$ws.Range("M77").Value = 54

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 738.1667
$ws.Range("I64").Value = 801
$ws.Range("J64").Value = 612.5
$ws.Range("K64").Value = 801
$ws.Range("L64").Value = 612.5
$ws.Range("M64").Value = -576
$ws.Range("N64").Value = -1062.5
$ws.Range("H67").Value = 738.1667
$ws.Range("I67").Value = 801
$ws.Range("J67").Value = 612.5
$ws.Range("K67").Value = 801
$ws.Range("L67").Value = 612.5
$ws.Range("M67").Value = -21
$ws.Range("N67").Value = -2172.5
$ws.Range("H76").Value = 11000
$ws.Range("J76").Value = 11000
$ws.Range("L76").Value = 11000
$ws.Range("N76").Value = -11630
$ws.Range("H79").Value = 11000
$ws.Range("J79").Value = 11000
$ws.Range("L79").Value = 11000
$ws.Range("N79").Value = -13184
$ws.Range("H94").Value = 1922.6428
$ws.Range("I94").Value = 1310.6364
$ws.Range("K94").Value = 1310.6364
$ws.Range("M94").Value = -859.6364000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 2900.7693
$ws.Range("I99").Value = 2655.4546
$ws.Range("K99").Value = 2655.4546
$ws.Range("M99").Value = -1157.4546
$ws.Range("H126").Value = 2900.7693
$ws.Range("I126").Value = 2655.4546
$ws.Range("K126").Value = 7966.3638
$ws.Range("M126").Value = -5496.3638

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 108.70588
$ws.Range("I6").Value = 108.70588
$ws.Range("K6").Value = 326.11764
$ws.Range("M6").Value = -213.11764
$ws.Range("H44").Value = 568.6667
$ws.Range("I44").Value = 568.6667
$ws.Range("K44").Value = 1706.0001
$ws.Range("M44").Value = -1308.0001
$ws.Range("H61").Value = 90
$ws.Range("I61").Value = 95
$ws.Range("J61").Value = 75
$ws.Range("K61").Value = 285
$ws.Range("L61").Value = 225
$ws.Range("M61").Value = -70
$ws.Range("N61").Value = -655
$ws.Range("H131").Value = 921.7143
$ws.Range("J131").Value = 967.5
$ws.Range("L131").Value = 2902.5
$ws.Range("N131").Value = -12982.5
$ws.Range("H132").Value = 4992.5
$ws.Range("I132").Value = 4985
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 44865
$ws.Range("L132").Value = 45000
$ws.Range("M132").Value = -42335
$ws.Range("N132").Value = -50060
$ws.Range("H137").Value = 2400
$ws.Range("I137").Value = 2400
$ws.Range("K137").Value = 7200
$ws.Range("M137").Value = -2100

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H31").Value = 2240
$ws.Range("I31").Value = 2200
$ws.Range("K31").Value = 2200
$ws.Range("M31").Value = -1908
$ws.Range("H37").Value = 2240
$ws.Range("I37").Value = 2200
$ws.Range("K37").Value = 2200
$ws.Range("M37").Value = -1923
$ws.Range("H97").Value = 3079.6667
$ws.Range("I97").Value = 2668.1428
$ws.Range("J97").Value = 3655.8
$ws.Range("K97").Value = 2668.1428
$ws.Range("L97").Value = 3655.8
$ws.Range("M97").Value = -2172.1428
$ws.Range("N97").Value = -4647.8
$ws.Range("H102").Value = 835.8095
$ws.Range("I102").Value = 835.8095
$ws.Range("K102").Value = 835.8095
$ws.Range("M102").Value = 786.1905
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 6814
$ws.Range("I16").Value = 6799.8335
$ws.Range("J16").Value = 6899
$ws.Range("K16").Value = 6799.8335
$ws.Range("L16").Value = 6899
$ws.Range("M16").Value = -6629.8335
$ws.Range("N16").Value = -7239
$ws.Range("H22").Value = 1365.6666
$ws.Range("I22").Value = 731.3333
$ws.Range("J22").Value = 2000
$ws.Range("K22").Value = 731.3333
$ws.Range("L22").Value = 2000
$ws.Range("M22").Value = -436.3333
$ws.Range("N22").Value = -2590
$ws.Range("H27").Value = 1365.6666
$ws.Range("I27").Value = 731.3333
$ws.Range("J27").Value = 2000
$ws.Range("K27").Value = 731.3333
$ws.Range("L27").Value = 2000
$ws.Range("M27").Value = -624.3333
$ws.Range("N27").Value = -2214
$ws.Range("H41").Value = 10000
$ws.Range("I41").Value = 10000
$ws.Range("K41").Value = 10000
$ws.Range("M41").Value = -9562
$ws.Range("H46").Value = 3855.8235
$ws.Range("I46").Value = 3322.6365
$ws.Range("J46").Value = 4833.3335
$ws.Range("K46").Value = 3322.6365
$ws.Range("L46").Value = 4833.3335
$ws.Range("M46").Value = -3134.6365
$ws.Range("N46").Value = -5209.3335
$ws.Range("H55").Value = 960.2
$ws.Range("I55").Value = 962.75
$ws.Range("K55").Value = 962.75
$ws.Range("M55").Value = -789.75
$ws.Range("H132").Value = 6751.375
$ws.Range("I132").Value = 3890
$ws.Range("K132").Value = 11670
$ws.Range("M132").Value = -9140

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3044.8
$ws.Range("I122").Value = 2883.111
$ws.Range("J122").Value = 4500
$ws.Range("K122").Value = 8649.332999999999
$ws.Range("L122").Value = 13500
$ws.Range("M122").Value = -6199.332999999999
$ws.Range("N122").Value = -18400
